$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.86"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.05"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.212"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05769"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.498"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.124"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8158"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8501"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1350"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06953"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03146"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02865"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09374"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.753"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001509"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04681"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005984"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006285"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001236"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004295"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006503"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.499"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3174"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1327"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03635"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006313"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1049"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003102"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007527"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005267"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.2902"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002336"

Write-Host "Applied all price/volume updates"